# Apply weekly shuffle of Fecha/Volumen/Precio values across rows 2-15.
# Only columns D, J, K, L, M, P change; each row's values come from a
# different (source) row in the original data - a fixed permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> source row (values to copy from)
$rowMap = @{
    2  = 4
    3  = 7
    4  = 10
    5  = 11
    6  = 15
    7  = 14
    8  = 8
    9  = 13
    10 = 5
    11 = 2
    12 = 9
    13 = 3
    14 = 6
    15 = 12
}

# Capture original values for D, J, K, L, M, P from every source row first,
# so that overwriting target rows does not clobber values still needed later.
$orig = @{}
foreach ($r in 2..15) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $vals = $orig[$source]

    $ws.Cells.Item($target, 4).Value2  = $vals.D   # D - Fecha
    $ws.Cells.Item($target, 10).Value2 = $vals.J   # J - Volumen
    $ws.Cells.Item($target, 11).Value2 = $vals.K   # K - Precio minimo
    $ws.Cells.Item($target, 12).Value2 = $vals.L   # L - Precio maximo
    $ws.Cells.Item($target, 13).Value2 = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($target, 16).Value2 = $vals.P   # P - Precio $/Kg
}
